# 完成最优潮流计算 - renumber bus indices (BUS1/BUS2) to 0-based numbering
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, 1),
    @(1, 2),
    @(2, 3),
    @(3, 4),
    @(4, 5),
    @(4, 6),
    @(5, 7),
    @(6, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("E13").Select()
